## Class02.pptx update:
##  1. Refresh the cached "datetimeFigureOut" date placeholder text
##     (12/11/18 -> 6/18/25) everywhere it appears: both slide masters,
##     every slide layout, the handout master and the notes master.
##  2. Two wording tweaks on slide 9's agenda textbox:
##       "present in Tuesday's class"  -> "present in next week's class"
##       "for both team members is equal" -> "for each team member is equal"

$p = $ppt.ActivePresentation

# --- helpers ---------------------------------------------------------

function Set-DatePlaceholderText($shapes, [string]$newText) {
    # ppPlaceholderDate = 16
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -eq -1 -and $shp.PlaceholderFormat.Type -eq 16) {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

function Replace-InShapeText($shape, [string]$oldText, [string]$newText) {
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    # Search using an ASCII-safe prefix of $oldText so the lookup isn't
    # tripped up by how curly punctuation round-trips through .Text.
    $anchorLen = [Math]::Min(30, $oldText.Length)
    $anchor = $oldText.Substring(0, $anchorLen)
    $idx = $full.IndexOf($anchor)
    if ($idx -lt 0) {
        throw "Could not find anchor text: $anchor"
    }
    $sub = $tr.Characters($idx + 1, $oldText.Length)
    $sub.Text = $newText
}

$newDate = "6/18/25"

# --- 1. date placeholders on both designs' slide masters + layouts ---

for ($di = 1; $di -le $p.Designs.Count; $di++) {
    $design = $p.Designs.Item($di)
    $master = $design.SlideMaster

    Set-DatePlaceholderText $master.Shapes $newDate

    for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
        $layout = $master.CustomLayouts.Item($li)
        Set-DatePlaceholderText $layout.Shapes $newDate
    }
}

# --- date placeholders on handout master + notes master --------------

Set-DatePlaceholderText $p.HandoutMaster.Shapes $newDate
Set-DatePlaceholderText $p.NotesMaster.Shapes $newDate

# --- 2. slide 9 wording tweaks ----------------------------------------

$slide9 = $p.Slides.Item(9)
$agenda = $slide9.Shapes.Item(2)   # "TextBox 5"

$old1 = "Every team of two students will present in Tuesday" + [char]0x2019 + "s class:"
$new1 = "Every team of two students will present in next week" + [char]0x2019 + "s class:"
Replace-InShapeText $agenda $old1 $new1

$old2 = "Make sure that the presentation time for both team members is equal"
$new2 = "Make sure that the presentation time for each team member is equal"
Replace-InShapeText $agenda $old2 $new2
